$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "VALOR MORA" amount (E11) ---
$ws.Range("E11").Value = 638867

# --- Update "Cant. Periodos" count (F13) ---
$ws.Range("F13").Value = 9

# --- Add a new period row (2509) to the worker detail table ---
# Insert a fresh blank row right after the current last data row (23), pushing
# the trailing blank/footer rows down by one.
$ws.Rows.Item(24).Insert()

# The row that used to be row 23 (period 2508, with the "closing" bottom-border
# style) now needs to move its formatting+values down into the new row 24.
$ws.Range("B23:J23").Copy($ws.Range("B24:J24"))

# Row 23 is no longer the last row in the table, so give it the regular
# "middle of table" formatting used by rows 16-22.
$ws.Range("B22:J22").Copy($ws.Range("B23:J23"))

# Row 23 keeps showing period 2508; row 24 becomes the new period 2509 entry.
$ws.Range("E23").Value = "2508"
$ws.Range("E24").Value = "2509"
$ws.Range("F24").Value = 74000
$ws.Range("G24").Value = 1850000
